$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.258144974708557
$ws.Range("B1").Value = 2.540671586990356
$ws.Range("C1").Value = 3.844839096069336
$ws.Range("D1").Value = 2.769141674041748
$ws.Range("E1").Value = 1.07072114944458
